# monitoramento.xlsx -- add projetos 4..12 rows + subetapa/natureza columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 5..13: duplicate formatting from row 4 (A:G), then set the values ---
$ws.Range("A4:G4").Copy($ws.Range("A5:G5")) | Out-Null
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "projeto 4"
$ws.Range("C5").Value = "categoria 4"
$ws.Range("D5").Value = 2001
$ws.Range("F5").Value = 43136

$ws.Range("A4:G4").Copy($ws.Range("A6:G6")) | Out-Null
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "projeto 5"
$ws.Range("C6").Value = "categoria 5"
$ws.Range("D6").Value = 2002
$ws.Range("F6").Value = 43137

$ws.Range("A4:G4").Copy($ws.Range("A7:G7")) | Out-Null
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "projeto 6"
$ws.Range("C7").Value = "categoria 6"
$ws.Range("D7").Value = 2003
$ws.Range("F7").Value = 43138

$ws.Range("A4:G4").Copy($ws.Range("A8:G8")) | Out-Null
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "projeto 7"
$ws.Range("C8").Value = "categoria 7"
$ws.Range("D8").Value = 2004
$ws.Range("F8").Value = 43139

$ws.Range("A4:G4").Copy($ws.Range("A9:G9")) | Out-Null
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "projeto 8"
$ws.Range("C9").Value = "categoria 8"
$ws.Range("D9").Value = 2005
$ws.Range("F9").Value = 43140

$ws.Range("A4:G4").Copy($ws.Range("A10:G10")) | Out-Null
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "projeto 9"
$ws.Range("C10").Value = "categoria 9"
$ws.Range("D10").Value = 2006
$ws.Range("F10").Value = 43141

$ws.Range("A4:G4").Copy($ws.Range("A11:G11")) | Out-Null
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "projeto 10"
$ws.Range("C11").Value = "categoria 10"
$ws.Range("D11").Value = 2007
$ws.Range("F11").Value = 43142

$ws.Range("A4:G4").Copy($ws.Range("A12:G12")) | Out-Null
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "projeto 11"
$ws.Range("C12").Value = "categoria 11"
$ws.Range("D12").Value = 2008
$ws.Range("F12").Value = 43143

$ws.Range("A4:G4").Copy($ws.Range("A13:G13")) | Out-Null
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "projeto 12"
$ws.Range("C13").Value = "categoria 12"
$ws.Range("D13").Value = 2009
$ws.Range("F13").Value = 43144

# --- new headers: I1 "subetapa", J1 "natureza" (copy alignment style from H1) ---
$ws.Range("H1").Copy($ws.Range("I1")) | Out-Null
$ws.Range("I1").Value = "subetapa"
$ws.Range("H1").Copy($ws.Range("J1")) | Out-Null
$ws.Range("J1").Value = "natureza"

# --- column I (subetapa number, underlined font) ---
$ws.Range("I2").Value = 1
$ws.Range("I2").Font.Underline = $true
$ws.Range("I3").Value = 1
$ws.Range("I3").Font.Underline = $true
$ws.Range("I4").Value = 1
$ws.Range("I4").Font.Underline = $true
$ws.Range("I5").Value = 1
$ws.Range("I5").Font.Underline = $true
$ws.Range("I6").Value = 1
$ws.Range("I6").Font.Underline = $true
$ws.Range("I7").Value = 2
$ws.Range("I7").Font.Underline = $true
$ws.Range("I8").Value = 2
$ws.Range("I8").Font.Underline = $true
$ws.Range("I9").Value = 2
$ws.Range("I9").Font.Underline = $true
$ws.Range("I10").Value = 9
$ws.Range("I10").Font.Underline = $true
$ws.Range("I11").Value = 10
$ws.Range("I11").Font.Underline = $true
$ws.Range("I12").Value = 11
$ws.Range("I12").Font.Underline = $true
$ws.Range("I13").Value = 12
$ws.Range("I13").Font.Underline = $true

# --- column J (natureza, publico/privado -- only rows 2..6 are populated) ---
$ws.Range("J2").Value = "publico"
$ws.Range("J3").Value = "privado"
$ws.Range("J4").Value = "publico"
$ws.Range("J5").Value = "privado"
$ws.Range("J6").Value = "publico"

# --- restore the previously-selected cell ---
$ws.Range("I9").Select() | Out-Null

Write-Host "edit complete"
